# Apply the "pack shipping orders" related update:
#  - Bins sheet gets 6 new bin-size rows (with unit dimensions encoded in the
#    bin name) describing standard box sizes used for packing.
#  - Items sheet header row is relabeled from Length/Width/Height/Weight to
#    Unit Length/Unit Width/Unit Height/Unit Weight.
#  - Selections on both sheets are updated to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$wsBins = $wb.Worksheets.Item("Bins")
$wsItems = $wb.Worksheets.Item("Items")

# --- Relabel header row on the "Items" sheet --------------------------------
# (done first so the new shared-string entries land in the same order as the
# target workbook: Unit Length/Width/Height/Weight before the bin sizes)
$wsItems.Range("D1").Value = "Unit Length"
$wsItems.Range("E1").Value = "Unit Width"
$wsItems.Range("F1").Value = "Unit Height"
$wsItems.Range("G1").Value = "Unit Weight"

# --- New bin rows on the "Bins" sheet ---------------------------------------
$newBins = @(
    @("20x20x30", 20, 20, 30, 50, 100),
    @("11x11x5", 11, 11, 5, 20, 10),
    @("22x22x18", 22, 22, 18, 20, 10),
    @("16x16x16", 16, 16, 16, 20, 10),
    @("24x24x24", 24, 24, 24, 20, 10),
    @("26x15x7", 26, 15, 7, 20, 10)
)

$row = 4
foreach ($bin in $newBins) {
    $wsBins.Cells.Item($row, 1).Value = $bin[0]
    $wsBins.Cells.Item($row, 2).Value = $bin[1]
    $wsBins.Cells.Item($row, 3).Value = $bin[2]
    $wsBins.Cells.Item($row, 4).Value = $bin[3]
    $wsBins.Cells.Item($row, 5).Value = $bin[4]
    $wsBins.Cells.Item($row, 6).Value = $bin[5]
    $row = $row + 1
}

# --- Restore the selections recorded in the workbook ------------------------
$wsBins.Activate()
$wsBins.Range("F10").Select()

$wsItems.Activate()
$wsItems.Range("D18").Select()
